$d = $word.ActiveDocument

# --- Locate the paragraph that ends in "repositorio" -----------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("repositorio", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'repositorio' paragraph"
}
$upTo = $d.Range(0, $searchRange.End)
$repoPara = $upTo.Paragraphs.Last

# --- Insert two new empty paragraphs right after it -------------------------
$repoPara.Range.InsertParagraphAfter()
$repoPara.Range.InsertParagraphAfter()

# --- Replace the last picture in the document with the new question text ---
$pic = $d.InlineShapes.Item($d.InlineShapes.Count)
$pic.Range.Text = "¿Cómo se ven los archivos individuales en la rama principal después de la fusión?"

# --- Remove the now-superfluous trailing empty paragraph -------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$markRange = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.Start)
$markRange.Delete()
